# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-20
$kValues = @{
    2  = 3
    3  = 3
    4  = 4
    5  = 5
    6  = 3
    7  = 2
    8  = 2
    9  = 6
    10 = 4
    11 = 3
    12 = 2
    13 = 3
    14 = 3
    15 = 2
    16 = 6
    17 = 5
    18 = 0
    19 = 2
    20 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
